# Payslip template: rename the merge-field placeholders in the "Standard Pay" /
# "Commission and Bonus" / "PAYE Tax" rows so they match the variable names
# used by the updated payroll-generation script (folder/output rework).
#
#   <Rate>    -> <Hourly Rate>   (C13, "Standard Pay" rate column)
#   <Bonuses> -> <Bonus>         (D14, "Commission and Bonus" row)
#   <Salary>  -> <salary>        (D13, "Standard Pay" current column)
#   <Taxes>   -> <tax>           (D21, "PAYE Tax" row)
#
# (Order matters: new shared-string entries are appended in first-use order,
# and the target workbook needs Hourly Rate, Bonus, salary, tax in that
# sequence.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C13").Value = "<Hourly Rate>"
$ws.Range("D14").Value = "<Bonus>"
$ws.Range("D13").Value = "<salary>"
$ws.Range("D21").Value = "<tax>"

# Match the author's final on-screen selection/scroll position.
$ws.Range("B22").Select()
